# Auto-generated script: apply scheduled market-data refresh to Chocobo_Profits sheets
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for the rows
# touched by this run, per-sheet, matching upstream scraped values.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 852.0707
$ws.Range("J129").Value = 903.06665
$ws.Range("L129").Value = 2709.19995
$ws.Range("N129").Value = -12709.19995

$ws.Range("H132").Value = 350076.56
$ws.Range("I132").Value = 5192.0386
$ws.Range("K132").Value = 15576.1158
$ws.Range("M132").Value = -13046.1158

$ws.Range("H137").Value = 1445190
$ws.Range("I137").Value = 2507850
$ws.Range("K137").Value = 7523550
$ws.Range("M137").Value = -7521000


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 902.1429000000001
$ws.Range("I2").Value = 656.17645
$ws.Range("J2").Value = 1947.5
$ws.Range("K2").Value = 656.17645
$ws.Range("L2").Value = 1947.5
$ws.Range("M2").Value = -543.17645
$ws.Range("N2").Value = -2173.5

$ws.Range("H32").Value = 4269.712
$ws.Range("I32").Value = 3683.074
$ws.Range("K32").Value = 3683.074
$ws.Range("M32").Value = -3396.074

$ws.Range("H45").Value = 3951.25
$ws.Range("I45").Value = 10011
$ws.Range("K45").Value = 10011
$ws.Range("M45").Value = -9634

$ws.Range("H61").Value = 1793
$ws.Range("I61").Value = 1793
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1793
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1581
$ws.Range("N61").ClearContents()

$ws.Range("H63").Value = 10658655
$ws.Range("I63").Value = 19789502
$ws.Range("K63").Value = 19789502
$ws.Range("M63").Value = -19788816

$ws.Range("H66").Value = 10658655
$ws.Range("I66").Value = 19789502
$ws.Range("K66").Value = 98947510
$ws.Range("M66").Value = -98944078

$ws.Range("H74").Value = 4348.567
$ws.Range("I74").Value = 4645.6665
$ws.Range("J74").Value = 3655.3333
$ws.Range("K74").Value = 4645.6665
$ws.Range("L74").Value = 3655.3333
$ws.Range("M74").Value = -3771.6665
$ws.Range("N74").Value = -5403.3333

$ws.Range("H77").Value = 4348.567
$ws.Range("I77").Value = 4645.6665
$ws.Range("J77").Value = 3655.3333
$ws.Range("K77").Value = 23228.3325
$ws.Range("L77").Value = 18276.6665
$ws.Range("M77").Value = -18860.3325
$ws.Range("N77").Value = -27012.6665

$ws.Range("H116").Value = 902.1429000000001
$ws.Range("I116").Value = 656.17645
$ws.Range("J116").Value = 1947.5
$ws.Range("K116").Value = 656.17645
$ws.Range("L116").Value = 1947.5
$ws.Range("M116").Value = 1637.82355
$ws.Range("N116").Value = -6535.5

$ws.Range("H132").Value = 2056.48
$ws.Range("I132").Value = 1220.65
$ws.Range("J132").Value = 5399.8
$ws.Range("K132").Value = 3661.95
$ws.Range("L132").Value = 16199.4
$ws.Range("M132").Value = -1131.95
$ws.Range("N132").Value = -21259.4

$ws.Range("H136").Value = 1793
$ws.Range("I136").Value = 1793
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5379
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2829
$ws.Range("N136").ClearContents()


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 902.1429000000001
$ws.Range("I3").Value = 656.17645
$ws.Range("J3").Value = 1947.5
$ws.Range("K3").Value = 656.17645
$ws.Range("L3").Value = 1947.5
$ws.Range("M3").Value = -542.17645
$ws.Range("N3").Value = -2175.5

$ws.Range("H38").Value = 14493.333
$ws.Range("J38").Value = 14493.333
$ws.Range("L38").Value = 14493.333
$ws.Range("N38").Value = -15325.333

$ws.Range("H86").Value = 2340.3333
$ws.Range("I86").Value = 2009.5454
$ws.Range("K86").Value = 2009.5454
$ws.Range("M86").Value = -886.5454

$ws.Range("H89").Value = 2340.3333
$ws.Range("I89").Value = 2009.5454
$ws.Range("K89").Value = 10047.727
$ws.Range("M89").Value = -4431.726999999999

$ws.Range("H105").Value = 5129823.5
$ws.Range("I105").Value = 5466099
$ws.Range("J105").Value = 1621
$ws.Range("K105").Value = 5466099
$ws.Range("L105").Value = 1621
$ws.Range("M105").Value = -5464352
$ws.Range("N105").Value = -5115

$ws.Range("H134").Value = 3093.1428
$ws.Range("I134").Value = 2377.4443
$ws.Range("J134").Value = 4381.4
$ws.Range("K134").Value = 7132.3329
$ws.Range("L134").Value = 13144.2
$ws.Range("M134").Value = -4597.3329
$ws.Range("N134").Value = -18214.2


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2406.3
$ws.Range("I31").Value = 1091.2084
$ws.Range("K31").Value = 1091.2084
$ws.Range("M31").Value = -796.2084

$ws.Range("H33").Value = 21345
$ws.Range("J33").Value = 40035
$ws.Range("L33").Value = 40035
$ws.Range("N33").Value = -40793

$ws.Range("H34").Value = 2406.3
$ws.Range("I34").Value = 1091.2084
$ws.Range("K34").Value = 1091.2084
$ws.Range("M34").Value = -889.2084

$ws.Range("H132").Value = 2499
$ws.Range("I132").Value = 1433.85
$ws.Range("K132").Value = 4301.549999999999
$ws.Range("M132").Value = -1771.549999999999


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 668203.9
$ws.Range("J5").Value = 1027700.75
$ws.Range("L5").Value = 3083102.25
$ws.Range("N5").Value = -3083326.25

$ws.Range("H12").Value = 107.13333
$ws.Range("J12").Value = 130.58333
$ws.Range("L12").Value = 391.74999
$ws.Range("N12").Value = -737.74999

$ws.Range("H107").Value = 84887.586
$ws.Range("I107").Value = 693.3333
$ws.Range("J107").Value = 112952.336
$ws.Range("K107").Value = 2079.9999
$ws.Range("L107").Value = 338857.008
$ws.Range("M107").Value = -159.9998999999998
$ws.Range("N107").Value = -342697.008

$ws.Range("H113").Value = 3907207
$ws.Range("I113").Value = 933.8421
$ws.Range("K113").Value = 2801.5263
$ws.Range("M113").Value = -631.5263

$ws.Range("H135").Value = 668203.9
$ws.Range("J135").Value = 1027700.75
$ws.Range("L135").Value = 9249306.75
$ws.Range("N135").Value = -9254376.75


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 2008.5
$ws.Range("I36").Value = 1017
$ws.Range("K36").Value = 1017
$ws.Range("M36").Value = -532

$ws.Range("H43").Value = 15345
$ws.Range("I43").Value = 1277.7778
$ws.Range("K43").Value = 1277.7778
$ws.Range("M43").Value = -1126.7778

$ws.Range("H46").Value = 30679.875
$ws.Range("J46").Value = 30679.875
$ws.Range("L46").Value = 30679.875
$ws.Range("N46").Value = -30991.875

$ws.Range("H80").Value = 35716788
$ws.Range("I80").Value = 50002344
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 50002344
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -50001346
$ws.Range("N80").Value = -4896

$ws.Range("H83").Value = 35716788
$ws.Range("I83").Value = 50002344
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 250011720
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -250006728
$ws.Range("N83").Value = -24484

$ws.Range("H132").Value = 15000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 45000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -50060


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4456.115
$ws.Range("I136").Value = 1952
$ws.Range("J136").Value = 6021.1875
$ws.Range("K136").Value = 5856
$ws.Range("L136").Value = 18063.5625
$ws.Range("M136").Value = -3306
$ws.Range("N136").Value = -23163.5625


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 19611036
$ws.Range("I132").Value = 1998
$ws.Range("K132").Value = 5994
$ws.Range("M132").Value = -3464

$ws.Range("H136").Value = 4298.7856
$ws.Range("I136").Value = 2518.7368
$ws.Range("J136").Value = 8056.6665
$ws.Range("K136").Value = 7556.2104
$ws.Range("L136").Value = 24169.9995
$ws.Range("M136").Value = -5006.2104
$ws.Range("N136").Value = -29269.9995

